# Auto-generated edit script applying the cryptos.xlsx price/volume update
# (commit: "Updated cryptos list on Thu Mar 28 03:43:57 UTC 2024 with GitHub Actions")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume) hold numeric-looking / percentage text
# (e.g. "1.00", "0.133", "  -1.82%  "). Force the range to Text format first
# so the Value setter does not silently coerce them into floating point
# numbers; the style is reset back to Normal afterwards so the saved cells
# stay styleless, matching the original inline-string cells.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "69.133.39"
$ws.Range("E2").Value = "  -1.82%  "
$ws.Range("D3").Value = "3.477.81"
$ws.Range("E3").Value = "  -3.66%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "577.43"
$ws.Range("E5").Value = "  -0.90%  "
$ws.Range("D6").Value = "181.13"
$ws.Range("E6").Value = "  -4.69%  "
$ws.Range("D7").Value = "3.467.64"
$ws.Range("E7").Value = "  -3.79%  "
$ws.Range("D8").Value = "0.607"
$ws.Range("E8").Value = "  -3.84%  "
$ws.Range("E9").Value = "  +0.04%  "
$ws.Range("D10").Value = "0.194"
$ws.Range("E10").Value = "  +5.35%  "
$ws.Range("D11").Value = "0.636"
$ws.Range("E11").Value = "  -3.97%  "
$ws.Range("D12").Value = "53.27"
$ws.Range("E12").Value = "  -5.08%  "
$ws.Range("D13").Value = "0.0000299"
$ws.Range("E13").Value = "  -4.06%  "
$ws.Range("D14").Value = "9.31"
$ws.Range("E14").Value = "  -4.09%  "
$ws.Range("D15").Value = "4.024.03"
$ws.Range("E15").Value = "  -4.19%  "
$ws.Range("D16").Value = "19.07"
$ws.Range("E16").Value = "  -4.19%  "
$ws.Range("D17").Value = "69.142.39"
$ws.Range("E17").Value = "  -1.82%  "
$ws.Range("D18").Value = "3.483.45"
$ws.Range("E18").Value = "  -3.70%  "
$ws.Range("D19").Value = "12.14"
$ws.Range("E19").Value = "  -4.16%  "
$ws.Range("E20").Value = "  -1.76%  "
$ws.Range("D21").Value = "533.25"
$ws.Range("E21").Value = "  +8.23%  "
$ws.Range("D22").Value = "0.999"
$ws.Range("E22").Value = "  -4.85%  "
$ws.Range("D23").Value = "18.33"
$ws.Range("E23").Value = "  -4.53%  "
$ws.Range("D24").Value = "4.45"
$ws.Range("E24").Value = "  +1.90%  "
$ws.Range("D25").Value = "4.81"
$ws.Range("E25").Value = "  -2.39%  "
$ws.Range("D26").Value = "95.42"
$ws.Range("E26").Value = "  -1.83%  "
$ws.Range("B27").Value = "RenderToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D27").Value = "10.93"
$ws.Range("E27").Value = "  -1.35%  "
$ws.Range("B28").Value = "ImmutableX"
$ws.Range("C28").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D28").Value = "2.91"
$ws.Range("E28").Value = "  -2.92%  "
$ws.Range("D29").Value = "8.99"
$ws.Range("E29").Value = "  -4.46%  "
$ws.Range("D30").Value = "31.69"
$ws.Range("E30").Value = "  -2.06%  "
$ws.Range("D31").Value = "7.13"
$ws.Range("E31").Value = "  -5.85%  "
$ws.Range("D32").Value = "12.32"
$ws.Range("E32").Value = "  +0.34%  "
$ws.Range("D33").Value = "63.42"
$ws.Range("E33").Value = "  -3.96%  "
$ws.Range("D35").Value = "529.40"
$ws.Range("E35").Value = "  -9.42%  "
$ws.Range("D36").Value = "0.400"
$ws.Range("E36").Value = "  -0.25%  "
$ws.Range("B37").Value = "Dai"
$ws.Range("C37").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D37").Value = "0.999"
$ws.Range("E37").Value = "  -0.27%  "
$ws.Range("B38").Value = "InjectiveProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D38").Value = "37.49"
$ws.Range("E38").Value = "  -4.17%  "
$ws.Range("D39").Value = "3.02"
$ws.Range("E39").Value = "  +2.65%  "
$ws.Range("D40").Value = "0.0₃0741"
$ws.Range("E40").Value = "  -9.15%  "
$ws.Range("D41").Value = "3.333.95"
$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").Value = "0.133"
$ws.Range("E42").Value = "  -3.09%  "
$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D43").Value = "3.34"
$ws.Range("E43").Value = "  -4.45%  "
$ws.Range("D44").Value = "3.51"
$ws.Range("E44").Value = "  +2.85%  "
$ws.Range("D45").Value = "2.94"
$ws.Range("E45").Value = "  -8.71%  "
$ws.Range("D46").Value = "2.91"
$ws.Range("E46").Value = "  -5.25%  "
$ws.Range("D47").Value = "0.0431"
$ws.Range("E47").Value = "  -3.27%  "
$ws.Range("E48").Value = "  -4.24%  "
$ws.Range("D49").Value = "8.89"
$ws.Range("E49").Value = "  -8.87%  "
$ws.Range("D50").Value = "1.00"
$ws.Range("E50").Value = "  +0.13%  "
$ws.Range("D51").Value = "136.34"
$ws.Range("E51").Value = "  -0.93%  "

$ws.Range("D2:E51").Style = "Normal"
